# Judging T5 enzyme data
# - Tidy up a couple of mis-entered cells on T0
# - Add PER (y/n) judgement column entries on T3
# - Fill in the full AG/AP/BG/BX/CBH/LAP/NAG/PPO/PER judgement table on T5

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# T0 : fix row 10 - CBH (F10) and LAP (G10) were entered in error,
# and NAG (H10) should be "x" (noisy) rather than "o" (substrate inhibition)
# ---------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("T0")
$ws0.Range("F10").ClearContents()
$ws0.Range("G10").ClearContents()
$ws0.Range("H10").Value = "x"

# ---------------------------------------------------------------
# T3 : add the PER (y/n) judgement column, plus a couple of missed
# PPO (I) judgements
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("T3")
$ws3.Range("J2").Value = "y"
$ws3.Range("J3").Value = "y"
$ws3.Range("I4").Value = "a"
$ws3.Range("J4").Value = "y"
$ws3.Range("J5").Value = "y"
$ws3.Range("J6").Value = "y"
$ws3.Range("J7").Value = "y"
$ws3.Range("I8").Value = "a"
$ws3.Range("J8").Value = "y"
$ws3.Range("I9").Value = "x"
$ws3.Range("J9").Value = "y"
$ws3.Range("I10").Value = "x"
$ws3.Range("J10").Value = "y"
$ws3.Range("I11").Value = "a"
$ws3.Range("J11").Value = "y"
$ws3.Range("I12").Value = "a"
$ws3.Range("J12").Value = "y"
$ws3.Range("I13").Value = "a"
$ws3.Range("J13").Value = "y"
$ws3.Range("I14").Value = "x"
$ws3.Range("J14").Value = "y"
$ws3.Range("J15").Value = "y"
$ws3.Range("I16").Value = "x"
$ws3.Range("J16").Value = "x"
$ws3.Range("I17").Value = "x"
$ws3.Range("J17").Value = "y"

# ---------------------------------------------------------------
# T5 : build out the whole judgement table (was completely empty)
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("T5")

# Pull over the header-row / sample-ID-column formatting used on the
# other two tabs so the new cells match (bold centered headers, and
# centered 12pt sample IDs down column A) without inventing new styles.
$ws3.Range("A1:J1").Copy()
$ws5.Range("A1:J1").PasteSpecial(-4122)
$ws3.Range("A2:A17").Copy()
$ws5.Range("A2:A17").PasteSpecial(-4122)

$ws5.Range("A1").Value = "ID"
$ws5.Range("B1").Value = "AG"
$ws5.Range("C1").Value = "AP"
$ws5.Range("D1").Value = "BG"
$ws5.Range("E1").Value = "BX"
$ws5.Range("F1").Value = "CBH"
$ws5.Range("G1").Value = "LAP"
$ws5.Range("H1").Value = "NAG"
$ws5.Range("I1").Value = "PPO"
$ws5.Range("J1").Value = "PER"

$ws5.Range("A2").Value = "4LXX"
$ws5.Range("C2").Value = "o"
$ws5.Range("D2").Value = "o"
$ws5.Range("F2").Value = "o"
$ws5.Range("G2").Value = "o"
$ws5.Range("H2").Value = "o"
$ws5.Range("I2").Value = "a"
$ws5.Range("J2").Value = "x"
$ws5.Range("K2").Value = "x = noisy, indicating low activity or bad data"

$ws5.Range("A3").Value = "5RRX"
$ws5.Range("D3").Value = "o"
$ws5.Range("I3").Value = "a"
$ws5.Range("J3").Value = "y"
$ws5.Range("K3").Value = "o = possible substrate inhibition"

$ws5.Range("A4").Value = "7LRX"
$ws5.Range("D4").Value = "o"
$ws5.Range("H4").Value = "o"
$ws5.Range("J4").Value = "y"
$ws5.Range("K4").Value = "a = salvageable with other errors"

$ws5.Range("A5").Value = "8LXX"
$ws5.Range("D5").Value = "o"
$ws5.Range("I5").Value = "a"
$ws5.Range("J5").Value = "y"
$ws5.Range("K5").Value = "y = generally negative activity, indicating no activity or need to improve methodology"

$ws5.Range("A6").Value = "14RRX"
$ws5.Range("D6").Value = "o"
$ws5.Range("H6").Value = "o"
$ws5.Range("I6").Value = "a"
$ws5.Range("J6").Value = "x"

$ws5.Range("A7").Value = "18RXX"
$ws5.Range("B7").Value = "o"
$ws5.Range("C7").Value = "o"
$ws5.Range("D7").Value = "o"
$ws5.Range("F7").Value = "o"
$ws5.Range("H7").Value = "o"
$ws5.Range("I7").Value = "a"
$ws5.Range("J7").Value = "y"

$ws5.Range("A8").Value = "20LRX"
$ws5.Range("H8").Value = "o"
$ws5.Range("I8").Value = "a"
$ws5.Range("J8").Value = "y"

$ws5.Range("A9").Value = "22LXX"
$ws5.Range("I9").Value = "a"
$ws5.Range("J9").Value = "x"

$ws5.Range("A10").Value = "25LRX"
$ws5.Range("B10").Value = "a"
$ws5.Range("F10").Value = "o"
$ws5.Range("G10").Value = "x"
$ws5.Range("H10").Value = "o"
$ws5.Range("I10").Value = "a"
$ws5.Range("J10").Value = "y"

$ws5.Range("A11").Value = "27RXX"
$ws5.Range("F11").Value = "o"
$ws5.Range("H11").Value = "o"
$ws5.Range("I11").Value = "a"
$ws5.Range("J11").Value = "y"

$ws5.Range("A12").Value = "32RXX"
$ws5.Range("B12").Value = "a"
$ws5.Range("G12").Value = "x"
$ws5.Range("H12").Value = "o"
$ws5.Range("I12").Value = "a"
$ws5.Range("J12").Value = "y"

$ws5.Range("A13").Value = "35LRX"
$ws5.Range("I13").Value = "a"
$ws5.Range("J13").Value = "y"

$ws5.Range("A14").Value = "45LRX"
$ws5.Range("D14").Value = "o"
$ws5.Range("H14").Value = "o"
$ws5.Range("I14").Value = "x"
$ws5.Range("J14").Value = "x"

$ws5.Range("A15").Value = "46RXX"
$ws5.Range("B15").Value = "x"
$ws5.Range("E15").Value = "x"
$ws5.Range("F15").Value = "o"
$ws5.Range("H15").Value = "o"

$ws5.Range("A16").Value = "47RRX"
$ws5.Range("B16").Value = "a"
$ws5.Range("H16").Value = "o"

$ws5.Range("A17").Value = "48LXX"
$ws5.Range("B17").Value = "a"
$ws5.Range("E17").Value = "a"
$ws5.Range("G17").Value = "x"
$ws5.Range("H17").Value = "o"

# ---------------------------------------------------------------
# Restore each tab's last-used selection, and finish with T5 active
# (it's the tab the author was working in when they saved).
# ---------------------------------------------------------------
$ws0.Range("B10").Select()
$ws3.Range("E23").Select()
$ws5.Range("I15").Select()
